{"js": "// Fix the typo \"Giver\" -> \"Given\" in the instruction paragraph\n// (\"Giver the above tables solve the following queries.\" ->\n//  \"Given the above tables solve the following queries.\").\nconst body = context.document.body;\n\nconst results = body.search(\"Giver\", { matchCase: true, matchWholeWord: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  // Replace the misspelled word in place; the surrounding run formatting\n  // (Calibri, size 24, color 000000, en-US) is preserved automatically.\n  results.items[0].insertText(\"Given\", \"Replace\");\n  await context.sync();\n}\n", "ps1": "# Fix the typo \"Giver\" -> \"Given\" in the instruction paragraph\n# (\"Giver the above tables solve the following queries.\" ->\n#  \"Given the above tables solve the following queries.\").\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.Text = \"Giver\"\n$find.Replacement.Text = \"Given\"\n$find.Execute($find.Text, $true, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n"}
